# TC01_CDS_Filter_Study-GECCO-OICR.xlsx : "Filter - Study - Test Suit"
#
# Semantic change: row 2 of the "startup" tab described the Cases tab
# ("CasesTab") but its query cell (B2) already held the Participants
# Cypher query - so the tab label is corrected to "ParticipantsTab".
# All other cells keep their original content. The active selection
# is moved from B3 to A2, and the row heights for the (now slightly
# taller, re-wrapped) rows 2-4 are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mislabeled tab name in A2: CasesTab -> ParticipantsTab
$ws.Range("A2").Value2 = "ParticipantsTab"

# Refresh the row heights for the wrapped-text rows (2-4) to match
# the re-wrapped text heights.
$ws.Rows.Item(2).RowHeight = 189
$ws.Rows.Item(3).RowHeight = 189
$ws.Rows.Item(4).RowHeight = 236.25

# Move the active selection to A2.
$ws.Range("A2").Select()
